$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily record was inserted at row 650 ("2026/01/18" / 日 / 19 / 18),
# pushing the existing rows 650..691 down to become rows 651..692.
$ws.Rows.Item(650).Insert()

# Populate the newly inserted row. Column A holds a date-look-alike string
# that must stay plain text (not be auto-converted to a date serial), so it
# is entered with a leading apostrophe; the cell style is then reset back to
# Normal so it matches the plain (unstyled) look of the surrounding data rows.
$ws.Cells.Item(650, 1).Value = "'2026/01/18"
$ws.Cells.Item(650, 1).Style = "Normal"
$ws.Cells.Item(650, 2).Value = "日"
$ws.Cells.Item(650, 3).Value = 19
$ws.Cells.Item(650, 4).Value = 18
